$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28-46 down to 29-47
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new record's data
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = 44651
$ws.Cells.Item(28, 4).NumberFormat = $ws.Cells.Item(29, 4).NumberFormat
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = 100112031
$ws.Cells.Item(28, 7).Value = "Poroto verde"
$ws.Cells.Item(28, 8).Value = "Magnum"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 140
$ws.Cells.Item(28, 11).Value = 20000
$ws.Cells.Item(28, 12).Value = 23000
$ws.Cells.Item(28, 13).Value = 21714
$ws.Cells.Item(28, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(28, 15).Value = "Región Metropolitana"
$ws.Cells.Item(28, 16).Value = 869
$ws.Cells.Item(28, 17).Value = 25
$ws.Cells.Item(28, 18).Value = "Hortaliza"
